# Auto-generated update of Leve profit figures across all profession sheets.
# Applies the recalculated currentAveragePrice / LevePrice / LeveProfit values
# produced by the scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 416.33334
$ws.Range("I12").Value = 375
$ws.Range("K12").Value = 375
$ws.Range("M12").Value = -205

$ws.Range("H39").Value = 1550.2632
$ws.Range("I39").Value = 860.375
$ws.Range("K39").Value = 2581.125
$ws.Range("M39").Value = -2285.125

$ws.Range("H55").Value = 917.1667
$ws.Range("I55").Value = 803
$ws.Range("K55").Value = 803
$ws.Range("M55").Value = -589

$ws.Range("H96").Value = 2999.2
$ws.Range("J96").Value = 2999.2
$ws.Range("L96").Value = 8997.599999999999
$ws.Range("N96").Value = -11743.6

$ws.Range("H132").Value = 7928.8
$ws.Range("I132").Value = 7374.591
$ws.Range("K132").Value = 22123.773
$ws.Range("M132").Value = -19593.773

$ws.Range("H135").Value = 1024.1111
$ws.Range("I135").Value = 1034.1428
$ws.Range("K135").Value = 9307.2852
$ws.Range("M135").Value = -6772.2852

$ws.Range("H138").Value = 2492.745
$ws.Range("I138").Value = 2020.0294
$ws.Range("K138").Value = 6060.0882
$ws.Range("M138").Value = -920.0882000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3075.875
$ws.Range("I45").Value = 1922
$ws.Range("K45").Value = 1922
$ws.Range("M45").Value = -1545

$ws.Range("H97").Value = 933.4
$ws.Range("I97").Value = 825.1429000000001
$ws.Range("K97").Value = 825.1429000000001
$ws.Range("M97").Value = -329.1429000000001

$ws.Range("H122").Value = 2337.25
$ws.Range("I122").Value = 1966.6666
$ws.Range("K122").Value = 5899.9998
$ws.Range("M122").Value = -3449.9998

$ws.Range("H132").Value = 37966.414
$ws.Range("I132").Value = 41693.46
$ws.Range("J132").Value = 5665.3335
$ws.Range("K132").Value = 125080.38
$ws.Range("L132").Value = 16996.0005
$ws.Range("M132").Value = -122550.38
$ws.Range("N132").Value = -22056.0005

$ws.Range("H135").Value = 49999.43
$ws.Range("J135").Value = 49999.43
$ws.Range("L135").Value = 49999.43
$ws.Range("N135").Value = -60139.43

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5973.3335
$ws.Range("I99").Value = 5400.476
$ws.Range("K99").Value = 5400.476
$ws.Range("M99").Value = -3902.476

$ws.Range("H107").Value = 1726.3182
$ws.Range("I107").Value = 1082.1666
$ws.Range("K107").Value = 1082.1666
$ws.Range("M107").Value = 837.8334

$ws.Range("H134").Value = 5807.091
$ws.Range("I134").Value = 5359.875
$ws.Range("J134").Value = 6999.6665
$ws.Range("K134").Value = 16079.625
$ws.Range("L134").Value = 20998.9995
$ws.Range("M134").Value = -13544.625
$ws.Range("N134").Value = -26068.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 598
$ws.Range("I22").Value = 412.625
$ws.Range("J22").Value = 809.8570999999999
$ws.Range("K22").Value = 412.625
$ws.Range("L22").Value = 809.8570999999999
$ws.Range("M22").Value = -62.625
$ws.Range("N22").Value = -1509.8571

$ws.Range("H58").Value = 79760.62
$ws.Range("J58").Value = 3876
$ws.Range("L58").Value = 3876
$ws.Range("N58").Value = -4282

$ws.Range("H86").Value = 5408.625
$ws.Range("I86").Value = 5636.6665
$ws.Range("J86").Value = 4724.5
$ws.Range("K86").Value = 5636.6665
$ws.Range("L86").Value = 4724.5
$ws.Range("M86").Value = -4513.6665
$ws.Range("N86").Value = -6970.5

$ws.Range("H89").Value = 5408.625
$ws.Range("I89").Value = 5636.6665
$ws.Range("J89").Value = 4724.5
$ws.Range("K89").Value = 28183.3325
$ws.Range("L89").Value = 23622.5
$ws.Range("M89").Value = -22567.3325
$ws.Range("N89").Value = -34854.5

$ws.Range("H94").Value = 1231
$ws.Range("I94").Value = 1000
$ws.Range("K94").Value = 1000
$ws.Range("M94").Value = -549

$ws.Range("H105").Value = 1313.7142
$ws.Range("I105").Value = 1032.6666
$ws.Range("K105").Value = 1032.6666
$ws.Range("M105").Value = 714.3334

$ws.Range("H107").Value = 2114.4
$ws.Range("I107").Value = 418.53333
$ws.Range("K107").Value = 418.53333
$ws.Range("M107").Value = 1501.46667

$ws.Range("H132").Value = 1554.1538
$ws.Range("I132").Value = 1567
$ws.Range("K132").Value = 4701
$ws.Range("M132").Value = -2171

$ws.Range("H134").Value = 72133.266
$ws.Range("I134").Value = 95280.73
$ws.Range("J134").Value = 8477.75
$ws.Range("K134").Value = 285842.19
$ws.Range("L134").Value = 25433.25
$ws.Range("M134").Value = -283307.19
$ws.Range("N134").Value = -30503.25

$ws.Range("H136").Value = 79760.62
$ws.Range("J136").Value = 3876
$ws.Range("L136").Value = 11628
$ws.Range("N136").Value = -16728

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 89.86667
$ws.Range("J12").Value = 86
$ws.Range("L12").Value = 258
$ws.Range("N12").Value = -604

$ws.Range("H131").Value = 12382.25
$ws.Range("J131").Value = 17086.857
$ws.Range("L131").Value = 51260.571
$ws.Range("N131").Value = -61340.571

$ws.Range("H132").Value = 1771
$ws.Range("J132").Value = 1894.6666
$ws.Range("L132").Value = 17051.9994
$ws.Range("N132").Value = -22111.9994

$ws.Range("H136").Value = 2198
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws.Range("H139").Value = 1538.5
$ws.Range("I139").Value = 1223.8182
$ws.Range("K139").Value = 3671.4546
$ws.Range("M139").Value = 1468.5454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4209.3
$ws.Range("I102").Value = 3548
$ws.Range("K102").Value = 3548
$ws.Range("M102").Value = -1926

$ws.Range("H107").Value = 72022.14
$ws.Range("I107").Value = 91028.17999999999
$ws.Range("K107").Value = 91028.17999999999
$ws.Range("M107").Value = -89108.17999999999

$ws.Range("H132").Value = 93246.37
$ws.Range("I132").Value = 102321
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 306963
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -304433
$ws.Range("N132").Value = -12560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4678.4287
$ws.Range("I7").Value = 2550
$ws.Range("J7").Value = 9999.5
$ws.Range("K7").Value = 2550
$ws.Range("L7").Value = 9999.5
$ws.Range("M7").Value = -2438
$ws.Range("N7").Value = -10223.5

$ws.Range("H55").Value = 1235.5
$ws.Range("I55").Value = 726.8570999999999
$ws.Range("J55").Value = 1744.1428
$ws.Range("K55").Value = 726.8570999999999
$ws.Range("L55").Value = 1744.1428
$ws.Range("M55").Value = -553.8570999999999
$ws.Range("N55").Value = -2090.1428

$ws.Range("H61").Value = 3124.92
$ws.Range("I61").Value = 2577.3809
$ws.Range("J61").Value = 5999.5
$ws.Range("K61").Value = 2577.3809
$ws.Range("L61").Value = 5999.5
$ws.Range("M61").Value = -2375.3809
$ws.Range("N61").Value = -6403.5

$ws.Range("H68").Value = 3336.0952
$ws.Range("I68").Value = 1874.6666
$ws.Range("K68").Value = 1874.6666
$ws.Range("M68").Value = -1125.6666

$ws.Range("H71").Value = 3336.0952
$ws.Range("I71").Value = 1874.6666
$ws.Range("K71").Value = 9373.333000000001
$ws.Range("M71").Value = -5629.333000000001

$ws.Range("H113").Value = 3124.92
$ws.Range("I113").Value = 2577.3809
$ws.Range("J113").Value = 5999.5
$ws.Range("K113").Value = 2577.3809
$ws.Range("L113").Value = 5999.5
$ws.Range("M113").Value = -407.3809000000001
$ws.Range("N113").Value = -10339.5

$ws.Range("H126").Value = 4678.4287
$ws.Range("I126").Value = 2550
$ws.Range("J126").Value = 9999.5
$ws.Range("K126").Value = 7650
$ws.Range("L126").Value = 29998.5
$ws.Range("M126").Value = -5180
$ws.Range("N126").Value = -34938.5

$ws.Range("H136").Value = 3888.7273
$ws.Range("I136").Value = 3027.6
$ws.Range("K136").Value = 9082.799999999999
$ws.Range("M136").Value = -6532.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1434.0555
$ws.Range("I107").Value = 1050.6
$ws.Range("K107").Value = 3151.8
$ws.Range("M107").Value = -1231.8

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
